{"js": "// Update the 25 division-expression cells in the worksheet's single table.\n// The table has 20 rows x 5 columns; only rows 0, 4, 8, 12, 16 (0-based)\n// contain text (\"a\u00f7b=\" style expressions) while the others are blank\n// spacer rows. We address each target cell by its fixed (row, col)\n// position and overwrite its text directly with the new expression -\n// this avoids any ambiguity from duplicate/overlapping old values.\nconst newValues = [\n  // row 0\n  [\"80\u00f79=\", \"42\u00f78=\", \"89\u00f72=\", \"27\u00f74=\", \"20\u00f74=\"],\n  // row 4\n  [\"78\u00f73=\", \"35\u00f75=\", \"20\u00f77=\", \"50\u00f74=\", \"32\u00f77=\"],\n  // row 8\n  [\"68\u00f72=\", \"48\u00f77=\", \"20\u00f73=\", \"32\u00f78=\", \"69\u00f74=\"],\n  // row 12\n  [\"50\u00f73=\", \"35\u00f74=\", \"43\u00f72=\", \"48\u00f72=\", \"84\u00f76=\"],\n  // row 16\n  [\"73\u00f78=\", \"21\u00f72=\", \"94\u00f77=\", \"65\u00f79=\", \"20\u00f73=\"],\n];\nconst rowIndices = [0, 4, 8, 12, 16];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nfor (let i = 0; i < rowIndices.length; i++) {\n  const r = rowIndices[i];\n  for (let c = 0; c < 5; c++) {\n    table.getCell(r, c).value = newValues[i][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the 25 division-expression cells in the document's single table.\n# The table has 20 rows x 5 columns; only rows 1, 5, 9, 13, 17 (1-based,\n# as COM indexes tables/cells) contain text (\"a\u00f7b=\" style expressions)\n# while the other rows are blank spacer rows. Each target cell is\n# addressed by its fixed (row, col) position and its text is overwritten\n# directly with the new expression - this avoids any ambiguity from\n# duplicate/overlapping old values during a text-search replace.\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n$newValues = @{\n    1  = @(\"80\u00f79=\", \"42\u00f78=\", \"89\u00f72=\", \"27\u00f74=\", \"20\u00f74=\")\n    5  = @(\"78\u00f73=\", \"35\u00f75=\", \"20\u00f77=\", \"50\u00f74=\", \"32\u00f77=\")\n    9  = @(\"68\u00f72=\", \"48\u00f77=\", \"20\u00f73=\", \"32\u00f78=\", \"69\u00f74=\")\n    13 = @(\"50\u00f73=\", \"35\u00f74=\", \"43\u00f72=\", \"48\u00f72=\", \"84\u00f76=\")\n    17 = @(\"73\u00f78=\", \"21\u00f72=\", \"94\u00f77=\", \"65\u00f79=\", \"20\u00f73=\")\n}\n\nforeach ($row in $newValues.Keys) {\n    $values = $newValues[$row]\n    for ($col = 1; $col -le 5; $col++) {\n        $t.Cell($row, $col).Range.Text = $values[$col - 1]\n    }\n}\n"}
